$d = $word.ActiveDocument

# Locate the target paragraph: the last paragraph in the document, which
# currently reads "(explicação da tabela de simbolos)" and contains the
# _GoBack bookmark.
$count = $d.Paragraphs.Count
$target = $d.Paragraphs.Item($count)
$rng = $target.Range

$rPr = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
$pPr = '<w:pPr><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>'

$paraA = '<w:p>' + $pPr + `
  '<w:r>' + $rPr + '<w:tab/><w:t xml:space="preserve">Quanto à implementação da tabela de símbolos, no que diz respeito à estrutura de dados foram implementadas duas estruturas em C, uma referente ao nó da tabela onde estão presentes o que vai ser impresso em cada linha da tabela, como o </w:t></w:r>' + `
  '<w:r>' + $rPr + '<w:t>Na</w:t></w:r>' + `
  '<w:r>' + $rPr + '<w:t>me, ParamTypes</w:t></w:r>' + `
  '<w:r>' + $rPr + '<w:t xml:space="preserve">, Type e param, a outra estrutura é referente à tabela em si, onde estão presentes as variáveis como o tipo, nome, array de parâmetros, e número de parâmetros, bem como um ponteiro para a estrutura do nó da tabela para </w:t></w:r>' + `
  '<w:r>' + $rPr + '<w:t>associar essa informação a um dado parâmetro ou método.</w:t></w:r>' + `
  '</w:p>'

$paraB = '<w:p>' + $pPr + `
  '<w:r>' + $rPr + '<w:tab/></w:r>' + `
  '<w:r>' + $rPr + '<w:t>Quanto aos algoritmos implementad</w:t></w:r>' + `
  '<w:r>' + $rPr + '<w:t>as diversas funções de verificação para os diversos símbolos terminais da gramática realizada no ficheiro yacc como o Program, FieldDecl, MethodDecl</w:t></w:r>' + `
  '<w:r>' + $rPr + '<w:t>,</w:t></w:r>' + `
  '<w:r>' + $rPr + '<w:t xml:space="preserve"> entre outros ou para a AST</w:t></w:r>' + `
  '<w:r>' + $rPr + '<w:t>, de modo a serem úteis para a implementação da</w:t></w:r>' + `
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
  '<w:r>' + $rPr + '<w:t xml:space="preserve"> AST anotada</w:t></w:r>' + `
  '<w:r>' + $rPr + '<w:t xml:space="preserve">. Para a realização da tabela de símbolos também foram implementadas funções de inserção e procura de elementos na tabela, assim como funções de inicialização, tanto da classe como </w:t></w:r>' + `
  '<w:r>' + $rPr + '<w:t>dos métodos presentes nesta e uma função para imprimir a tabela de símbolos da forma que é referido no enunciado.</w:t></w:r>' + `
  '</w:p>'

$xml = '<?xml version="1.0" standalone="yes"?>' + `
  '<?mso-application progid="Word.Document"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData>' + `
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:body>' + $paraA + $paraB + '</w:body>' + `
  '</w:document>' + `
  '</pkg:xmlData>' + `
  '</pkg:part>' + `
  '</pkg:package>'

$rng.InsertXML($xml)

# InsertXML placed the two new paragraphs just before the original
# paragraph mark; that original (now empty) paragraph mark survives as a
# trailing empty paragraph. Remove it by deleting the range spanning from
# the end of our new last paragraph up to (and including) that empty
# paragraph mark, which merges everything back into a single final
# paragraph carrying the original paragraph-mark identity/properties.
$newCount = $d.Paragraphs.Count
$lastReal = $d.Paragraphs.Item($newCount - 1)
$trailingEmpty = $d.Paragraphs.Item($newCount)
$mergeRange = $d.Range($lastReal.Range.End - 1, $trailingEmpty.Range.End - 1)
$mergeRange.Delete()

Write-Host "Final paragraph count:" $d.Paragraphs.Count
